$d = $word.ActiveDocument

function Find-Bookmark-Para {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $xml = $p.Range.WordOpenXML
        if ($xml -match "w:name=`"_GoBack`"") {
            return $i
        }
    }
    return -1
}

$bmIdx = Find-Bookmark-Para
if ($bmIdx -eq -1) {
    throw "Could not locate the paragraph containing the _GoBack bookmark"
}

# --- Insert the new paragraphs that go BEFORE the "Example:" / bookmark paragraph ---
# 1) a blank paragraph right after the "    ?>" paragraph
$bm = $d.Paragraphs.Item($bmIdx)
$bm.Range.InsertParagraphBefore()
$bmIdx = Find-Bookmark-Para

# 2) "For Loops in PHP" paragraph
$bm = $d.Paragraphs.Item($bmIdx)
$bm.Range.InsertBefore("For Loops in PHP`r")
$bmIdx = Find-Bookmark-Para

# 3) "We can use it to repeat ..." paragraph
$bm = $d.Paragraphs.Item($bmIdx)
$bm.Range.InsertBefore("We can use it to repeat a series of instructions instead of typing a lot of print of echo many times.`r")
$bmIdx = Find-Bookmark-Para

# 4) add "Example:" run to the start of the bookmark paragraph itself (same paragraph)
$bm = $d.Paragraphs.Item($bmIdx)
$bm.Range.InsertBefore("Example:")
$bmIdx = Find-Bookmark-Para

# --- Insert the new paragraphs that go AFTER the bookmark paragraph ---
$bm = $d.Paragraphs.Item($bmIdx)
for ($k = 0; $k -lt 7; $k++) {
    $bm.Range.InsertParagraphAfter()
}

$p1 = $d.Paragraphs.Item($bmIdx + 1)
$p1.Range.Text = "This will print a list of leap years"

$xmlNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$phpFrag = "<pkg:package xmlns:pkg=`"http://schemas.microsoft.com/office/2006/xmlPackage`"><pkg:part pkg:name=`"/word/document.xml`" pkg:contentType=`"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml`"><pkg:xmlData><w:document $xmlNs><w:body><w:p><w:proofErr w:type=`"gramStart`"/><w:r><w:t>&lt;?php</w:t></w:r><w:proofErr w:type=`"gramEnd`"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
$p2 = $d.Paragraphs.Item($bmIdx + 2)
$p2.Range.InsertXML($phpFrag)

$p3 = $d.Paragraphs.Item($bmIdx + 3)
$p3.Range.Text = "      for (`$leap = 2004; `$leap < 2050; `$leap = `$leap + 4) {"

$p4 = $d.Paragraphs.Item($bmIdx + 4)
$p4.Range.Text = "        echo `"<p>`$leap</p>`";"

$p5 = $d.Paragraphs.Item($bmIdx + 5)
$p5.Range.Text = "      }"

$p6 = $d.Paragraphs.Item($bmIdx + 6)
$p6.Range.Text = "    ?>"

# $bmIdx + 7 stays the trailing blank paragraph

Write-Output "Done. bmIdx=$bmIdx ParagraphCount=$($d.Paragraphs.Count)"
